{"js": "// Fix the typo in the document title: \"Task 2\" -> \"Task 1\"\n// (commit message: \"Edit typo in title\")\nconst body = context.document.body;\nconst results = body.search(\"Task 2: Dead Function Analysis\", { matchCase: true });\nresults.load(\"items,text\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"Task 1: Dead Function Analysis\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Fix the typo in the document title: \"Task 2\" -> \"Task 1\"\n# (commit message: \"Edit typo in title\")\n$d = $word.ActiveDocument\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"Task 2: Dead Function Analysis\"\n$find.Replacement.Text = \"Task 1: Dead Function Analysis\"\n$find.Execute([ref]\"Task 2: Dead Function Analysis\", $false, $false, $false, $false, $false, $true, 1, $false, \"Task 1: Dead Function Analysis\", 2)\n"}
